$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5498.1816
$ws.Range("I32").Value = 4830.8335
$ws.Range("K32").Value = 4830.8335
$ws.Range("M32").Value = -4504.8335

$ws.Range("H98").Value = 638.6579
$ws.Range("I98").Value = 628.8919
$ws.Range("K98").Value = 628.8919
$ws.Range("M98").Value = 869.1081

$ws.Range("H100").Value = 1509.1765
$ws.Range("J100").Value = 2356
$ws.Range("L100").Value = 2356
$ws.Range("N100").Value = -3438

$ws.Range("H111").Value = 2313.0527
$ws.Range("I111").Value = 592.5
$ws.Range("J111").Value = 2515.4707
$ws.Range("K111").Value = 1777.5
$ws.Range("L111").Value = 7546.4121
$ws.Range("M111").Value = 1289.5
$ws.Range("N111").Value = -13680.4121

$ws.Range("H122").Value = 638.6579
$ws.Range("I122").Value = 628.8919
$ws.Range("K122").Value = 1886.6757
$ws.Range("M122").Value = 563.3243000000002

$ws.Range("H132").Value = 5445.577
$ws.Range("I132").Value = 1445.4445
$ws.Range("K132").Value = 4336.333500000001
$ws.Range("M132").Value = -1806.333500000001

$ws.Range("H134").Value = 84263.5
$ws.Range("J134").Value = 84263.5
$ws.Range("L134").Value = 84263.5
$ws.Range("N134").Value = -94403.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3471.9167
$ws.Range("I32").Value = 3333.0212
$ws.Range("K32").Value = 3333.0212
$ws.Range("M32").Value = -3046.0212

$ws.Range("H41").Value = 22587.5
$ws.Range("I41").Value = 350
$ws.Range("K41").Value = 350
$ws.Range("M41").Value = 64

$ws.Range("H45").Value = 204497.7
$ws.Range("I45").Value = 402195.6
$ws.Range("K45").Value = 402195.6
$ws.Range("M45").Value = -401818.6

$ws.Range("H61").Value = 9080.482
$ws.Range("J61").Value = 9369.799999999999
$ws.Range("L61").Value = 9369.799999999999
$ws.Range("N61").Value = -9793.799999999999

$ws.Range("H74").Value = 3313.0652
$ws.Range("I74").Value = 2307.9211
$ws.Range("J74").Value = 8087.5
$ws.Range("K74").Value = 2307.9211
$ws.Range("L74").Value = 8087.5
$ws.Range("M74").Value = -1433.9211
$ws.Range("N74").Value = -9835.5

$ws.Range("H77").Value = 3313.0652
$ws.Range("I77").Value = 2307.9211
$ws.Range("J77").Value = 8087.5
$ws.Range("K77").Value = 11539.6055
$ws.Range("L77").Value = 40437.5
$ws.Range("M77").Value = -7171.6055
$ws.Range("N77").Value = -49173.5

$ws.Range("H136").Value = 9080.482
$ws.Range("J136").Value = 9369.799999999999
$ws.Range("L136").Value = 28109.4
$ws.Range("N136").Value = -33209.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3999.1428
$ws.Range("I7").Value = 3997
$ws.Range("K7").Value = 3997
$ws.Range("M7").Value = -3884

$ws.Range("H86").Value = 2819.875
$ws.Range("I86").Value = 2306.3076
$ws.Range("K86").Value = 2306.3076
$ws.Range("M86").Value = -1183.3076

$ws.Range("H89").Value = 2819.875
$ws.Range("I89").Value = 2306.3076
$ws.Range("K89").Value = 11531.538
$ws.Range("M89").Value = -5915.538

$ws.Range("H134").Value = 3182.224
$ws.Range("I134").Value = 3223.1404
$ws.Range("J134").Value = 850
$ws.Range("K134").Value = 9669.421200000001
$ws.Range("L134").Value = 2550
$ws.Range("M134").Value = -7134.421200000001
$ws.Range("N134").Value = -7620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 6624.75
$ws.Range("I12").Value = 6624.75
$ws.Range("K12").Value = 6624.75
$ws.Range("M12").Value = -6454.75

$ws.Range("H16").Value = 4249.9375
$ws.Range("I16").Value = 2263.75
$ws.Range("K16").Value = 2263.75
$ws.Range("M16").Value = -1976.75

$ws.Range("H31").Value = 6829.4287
$ws.Range("I31").Value = 6199.8335
$ws.Range("K31").Value = 6199.8335
$ws.Range("M31").Value = -5904.8335

$ws.Range("H34").Value = 6829.4287
$ws.Range("I34").Value = 6199.8335
$ws.Range("K34").Value = 6199.8335
$ws.Range("M34").Value = -5997.8335

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null

$ws.Range("H113").Value = 4249.9375
$ws.Range("I113").Value = 2263.75
$ws.Range("K113").Value = 2263.75
$ws.Range("M113").Value = -93.75

$ws.Range("H132").Value = 3500.875
$ws.Range("I132").Value = 2851.4546
$ws.Range("J132").Value = 4929.6
$ws.Range("K132").Value = 8554.363799999999
$ws.Range("L132").Value = 14788.8
$ws.Range("M132").Value = -6024.363799999999
$ws.Range("N132").Value = -19848.8

$ws.Range("H141").Value = 31649
$ws.Range("J141").Value = 30200
$ws.Range("L141").Value = 30200
$ws.Range("N141").Value = -40560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1821.9524
$ws.Range("I5").Value = 1905.3334
$ws.Range("K5").Value = 5716.0002
$ws.Range("M5").Value = -5604.0002

$ws.Range("H38").Value = 360.41666
$ws.Range("J38").Value = 393.0909
$ws.Range("L38").Value = 1179.2727
$ws.Range("N38").Value = -1873.2727

$ws.Range("H55").Value = 701.3333
$ws.Range("J55").Value = 768.3333
$ws.Range("L55").Value = 2304.9999
$ws.Range("N55").Value = -2658.9999

$ws.Range("H103").Value = 778.8
$ws.Range("I103").Value = 631.6667
$ws.Range("K103").Value = 1895.0001
$ws.Range("M103").Value = -1016.0001

$ws.Range("H106").Value = 10915.667
$ws.Range("J106").Value = 3500
$ws.Range("L106").Value = 10500
$ws.Range("N106").Value = -12392

$ws.Range("H129").Value = 926804.25
$ws.Range("I129").Value = 202330.2
$ws.Range("K129").Value = 606990.6000000001
$ws.Range("M129").Value = -601990.6000000001

$ws.Range("H132").Value = 1922.6957
$ws.Range("I132").Value = 1656.5454
$ws.Range("K132").Value = 14908.9086
$ws.Range("M132").Value = -12378.9086

$ws.Range("H135").Value = 1821.9524
$ws.Range("I135").Value = 1905.3334
$ws.Range("K135").Value = 17148.0006
$ws.Range("M135").Value = -14613.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 32251
$ws.Range("J96").Value = 32251
$ws.Range("L96").Value = 32251
$ws.Range("N96").Value = -37743

$ws.Range("H102").Value = 3894.5454
$ws.Range("I102").Value = 3292.5
$ws.Range("K102").Value = 3292.5
$ws.Range("M102").Value = -1670.5

$ws.Range("H113").Value = 287832.94
$ws.Range("I113").Value = 446419.44
$ws.Range("K113").Value = 446419.44
$ws.Range("M113").Value = -444249.44

$ws.Range("H132").Value = 3574.7144
$ws.Range("I132").Value = 1604.6
$ws.Range("K132").Value = 4813.799999999999
$ws.Range("M132").Value = -2283.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 7089
$ws.Range("I39").Value = 7089
$ws.Range("K39").Value = 7089
$ws.Range("M39").Value = -6629

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null

$ws.Range("H82").Value = 1023.5833
$ws.Range("J82").Value = 1666
$ws.Range("L82").Value = 1666
$ws.Range("N82").Value = -2388

$ws.Range("H85").Value = 1023.5833
$ws.Range("J85").Value = 1666
$ws.Range("L85").Value = 1666
$ws.Range("N85").Value = -4162

$ws.Range("H132").Value = 10632.514
$ws.Range("I132").Value = 10814.6875
$ws.Range("J132").Value = 9466.6
$ws.Range("K132").Value = 32444.0625
$ws.Range("L132").Value = 28399.8
$ws.Range("M132").Value = -29914.0625
$ws.Range("N132").Value = -33459.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2503.375
$ws.Range("J96").Value = 2687.25
$ws.Range("L96").Value = 2687.25
$ws.Range("N96").Value = -5433.25
